$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 186 and 187 (match ids 7648957 / 7648958) had their data swapped
#    in the source feed (everything except the sequential "id" in column A,
#    which stays tied to the row position). Swap columns B:AC between row
#    186 and row 187.
# ---------------------------------------------------------------------------
$lastCol = 29  # column AC
for ($c = 2; $c -le $lastCol; $c++) {
    $v186 = $ws.Cells.Item(186, $c).Value2
    $v187 = $ws.Cells.Item(187, $c).Value2
    $ws.Cells.Item(186, $c).Value = $v187
    $ws.Cells.Item(187, $c).Value = $v186
}

# ---------------------------------------------------------------------------
# 2) Row 225 (id 223, matchId 7641726, Universidad Guadalajara vs
#    Tepatitlan FC) is replaced by a new upcoming fixture (matchId 7641723,
#    Cancun FC vs Dorados). Three more upcoming fixtures are appended as new
#    rows 226-228, the last of which (228) is the original 7641726 fixture
#    with refreshed (closing) odds.
# ---------------------------------------------------------------------------

# New row 225: Cancun FC vs Dorados
$ws.Range("A225").Value = 223
$ws.Range("B225").Value = 7641723
$ws.Range("C225").Value = "Mexico Liga de Expansion"
$ws.Range("D225").Value = "Mexico Liga de Expansion"
$ws.Range("E225").Value = 45391.92013888889
$ws.Range("F225").Value = "Cancun FC"
$ws.Range("G225").Value = "Dorados"
$ws.Range("K225").Value = 1.4
$ws.Range("L225").Value = 4
$ws.Range("M225").Value = 6.5
$ws.Range("N225").Value = 1.3
$ws.Range("O225").Value = 5.25
$ws.Range("P225").Value = 9
$ws.Range("Q225").Value = -1.75
$ws.Range("R225").Value = 2
$ws.Range("S225").Value = 1.8
$ws.Range("T225").Value = 3
$ws.Range("U225").Value = 1.85
$ws.Range("V225").Value = 1.95
$ws.Range("W225").Value = 0
$ws.Range("X225").Value = 0
$ws.Range("Y225").Value = 0
$ws.Range("Z225").Value = 0
$ws.Range("AA225").Value = 0

# New row 226: Club Atletico La Paz vs Venados FC
$ws.Range("A226").Value = 224
$ws.Range("B226").Value = 7641724
$ws.Range("C226").Value = "Mexico Liga de Expansion"
$ws.Range("D226").Value = "Mexico Liga de Expansion"
$ws.Range("E226").Value = 45392.00347222222
$ws.Range("F226").Value = "Club Atletico La Paz"
$ws.Range("G226").Value = "Venados FC"
$ws.Range("K226").Value = 2.3
$ws.Range("L226").Value = 3.25
$ws.Range("M226").Value = 2.7
$ws.Range("N226").Value = 2.9
$ws.Range("O226").Value = 3.4
$ws.Range("P226").Value = 2.375
$ws.Range("Q226").Value = 0.25
$ws.Range("R226").Value = 1.775
$ws.Range("S226").Value = 2.025
$ws.Range("T226").Value = 2.5
$ws.Range("U226").Value = 1.925
$ws.Range("V226").Value = 1.875
$ws.Range("W226").Value = 0
$ws.Range("X226").Value = 0
$ws.Range("Y226").Value = 0
$ws.Range("Z226").Value = 0
$ws.Range("AA226").Value = 0

# New row 227: Tlaxcala FC vs Mineros de Zacatecas
$ws.Range("A227").Value = 225
$ws.Range("B227").Value = 7641725
$ws.Range("C227").Value = "Mexico Liga de Expansion"
$ws.Range("D227").Value = "Mexico Liga de Expansion"
$ws.Range("E227").Value = 45392.92013888889
$ws.Range("F227").Value = "Tlaxcala FC"
$ws.Range("G227").Value = "Mineros de Zacatecas"
$ws.Range("K227").Value = 2.75
$ws.Range("L227").Value = 3.25
$ws.Range("M227").Value = 2.25
$ws.Range("N227").Value = 3.6
$ws.Range("O227").Value = 3.6
$ws.Range("P227").Value = 1.95
$ws.Range("Q227").Value = 0.5
$ws.Range("R227").Value = 1.875
$ws.Range("S227").Value = 1.925
$ws.Range("T227").Value = 2.75
$ws.Range("U227").Value = 1.9
$ws.Range("V227").Value = 1.9
$ws.Range("W227").Value = 0
$ws.Range("X227").Value = 0
$ws.Range("Y227").Value = 0
$ws.Range("Z227").Value = 0
$ws.Range("AA227").Value = 0

# New row 228: Universidad Guadalajara vs Tepatitlan FC (the old row 225
# fixture, now re-priced with refreshed odds)
$ws.Range("A228").Value = 226
$ws.Range("B228").Value = 7641726
$ws.Range("C228").Value = "Mexico Liga de Expansion"
$ws.Range("D228").Value = "Mexico Liga de Expansion"
$ws.Range("E228").Value = 45393.00347222222
$ws.Range("F228").Value = "Universidad Guadalajara"
$ws.Range("G228").Value = "Tepatitlan FC"
$ws.Range("K228").Value = 1.25
$ws.Range("L228").Value = 5.5
$ws.Range("M228").Value = 7.5
$ws.Range("N228").Value = 1.3
$ws.Range("O228").Value = 5.5
$ws.Range("P228").Value = 8
$ws.Range("Q228").Value = -1.5
$ws.Range("R228").Value = 1.825
$ws.Range("S228").Value = 1.975
$ws.Range("T228").Value = 2.75
$ws.Range("U228").Value = 1.825
$ws.Range("V228").Value = 1.975
$ws.Range("W228").Value = 0
$ws.Range("X228").Value = 0
$ws.Range("Y228").Value = 0
$ws.Range("Z228").Value = 0
$ws.Range("AA228").Value = 0

# ---------------------------------------------------------------------------
# Copy the formatting (styles) of row 225's A/E cells (bold/border/center
# style for id column, custom date format for Date column) onto the three
# newly appended rows, matching every other data row in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A225").Copy()
$ws.Range("A226:A228").PasteSpecial(-4122)
$ws.Range("E225").Copy()
$ws.Range("E226:E228").PasteSpecial(-4122)
$excel.CutCopyMode = 0
